$d = $word.ActiveDocument

$replacements = @(
    @{old="132÷5=26, 2"; new="427÷7=61, 0"},
    @{old="247÷8=30, 7"; new="128÷8=16, 0"},
    @{old="693÷9=77, 0"; new="610÷4=152, 2"},
    @{old="868÷6=144, 4"; new="944÷9=104, 8"},
    @{old="244÷6=40, 4"; new="135÷7=19, 2"},
    @{old="666÷4=166, 2"; new="372÷5=74, 2"},
    @{old="507÷6=84, 3"; new="689÷6=114, 5"},
    @{old="772÷9=85, 7"; new="448÷5=89, 3"},
    @{old="515÷8=64, 3"; new="318÷7=45, 3"},
    @{old="705÷5=141, 0"; new="875÷7=125, 0"},
    @{old="123÷8=15, 3"; new="641÷6=106, 5"},
    @{old="491÷3=163, 2"; new="231÷7=33, 0"},
    @{old="690÷4=172, 2"; new="238÷3=79, 1"},
    @{old="468÷5=93, 3"; new="526÷7=75, 1"},
    @{old="629÷8=78, 5"; new="616÷7=88, 0"},
    @{old="423÷5=84, 3"; new="817÷6=136, 1"},
    @{old="695÷2=347, 1"; new="315÷7=45, 0"},
    @{old="193÷3=64, 1"; new="146÷2=73, 0"},
    @{old="420÷7=60, 0"; new="459÷9=51, 0"},
    @{old="328÷3=109, 1"; new="967÷8=120, 7"},
    @{old="387÷3=129, 0"; new="292÷3=97, 1"},
    @{old="153÷8=19, 1"; new="826÷2=413, 0"},
    @{old="472÷8=59, 0"; new="306÷3=102, 0"},
    @{old="549÷4=137, 1"; new="667÷4=166, 3"},
    @{old="439÷7=62, 5"; new="635÷8=79, 3"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

Write-Host "Done"
